$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: id=1, name=Нори
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Нори"
$ws.Cells.Item(2, 3).Value = 10000
$ws.Cells.Item(2, 4).Value = "г"
$ws.Cells.Item(2, 5).Value = 1.1

# Row 3: id=2, name=Рис
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Рис"
$ws.Cells.Item(3, 3).Value = 10000
$ws.Cells.Item(3, 4).Value = "г"
$ws.Cells.Item(3, 5).Value = 0.18

# Row 4: id=3, name=банан
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "банан"
$ws.Cells.Item(4, 3).Value = 10000
$ws.Cells.Item(4, 4).Value = "г"
$ws.Cells.Item(4, 5).Value = 0.18

# Row 5: id=4, name=васаби
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "васаби"
$ws.Cells.Item(5, 3).Value = 10000
$ws.Cells.Item(5, 4).Value = "г"
$ws.Cells.Item(5, 5).Value = 0.5

# Row 6: id=5, name=вода
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "вода"
$ws.Cells.Item(6, 3).Value = 10000
$ws.Cells.Item(6, 4).Value = "г"
$ws.Cells.Item(6, 5).Value = 0

# Row 7: id=6, name=икра масага
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "икра масага"
$ws.Cells.Item(7, 3).Value = 10000
$ws.Cells.Item(7, 4).Value = "г"
$ws.Cells.Item(7, 5).Value = 1.33

# Row 8: id=7, name=имбир
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "имбир"
$ws.Cells.Item(8, 3).Value = 10000
$ws.Cells.Item(8, 4).Value = "г"
$ws.Cells.Item(8, 5).Value = 0.24

# Row 9: id=8, name=имбирь
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "имбирь"
$ws.Cells.Item(9, 3).Value = 10000
$ws.Cells.Item(9, 4).Value = "г"
$ws.Cells.Item(9, 5).Value = 0.24

# Row 10: id=9, name=кап семка
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "кап семка"
$ws.Cells.Item(10, 3).Value = 10000
$ws.Cells.Item(10, 4).Value = "г"
$ws.Cells.Item(10, 5).Value = 2.69

# Row 11: id=10, name=капч лосось
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "капч лосось"
$ws.Cells.Item(11, 3).Value = 10000
$ws.Cells.Item(11, 4).Value = "г"
$ws.Cells.Item(11, 5).Value = 0.49

# Row 12: id=11, name=капчен кур
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "капчен кур"
$ws.Cells.Item(12, 3).Value = 10000
$ws.Cells.Item(12, 4).Value = "г"
$ws.Cells.Item(12, 5).Value = 0.49

# Row 13: id=12, name=киви
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "киви"
$ws.Cells.Item(13, 3).Value = 10000
$ws.Cells.Item(13, 4).Value = "г"
$ws.Cells.Item(13, 5).Value = 0.2

# Row 14: id=13, name=клубника
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "клубника"
$ws.Cells.Item(14, 3).Value = 10000
$ws.Cells.Item(14, 4).Value = "г"
$ws.Cells.Item(14, 5).Value = 0.25

# Row 15: id=14, name=кляр
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "кляр"
$ws.Cells.Item(15, 3).Value = 10000
$ws.Cells.Item(15, 4).Value = "г"
$ws.Cells.Item(15, 5).Value = 0.12

# Row 16: id=15, name=краб
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "краб"
$ws.Cells.Item(16, 3).Value = 10000
$ws.Cells.Item(16, 4).Value = "г"
$ws.Cells.Item(16, 5).Value = 1.07

# Row 17: id=16, name=кунжут
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "кунжут"
$ws.Cells.Item(17, 3).Value = 10000
$ws.Cells.Item(17, 4).Value = "г"
$ws.Cells.Item(17, 5).Value = 0.25

# Row 18: id=17, name=курица
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "курица"
$ws.Cells.Item(18, 3).Value = 10000
$ws.Cells.Item(18, 4).Value = "г"
$ws.Cells.Item(18, 5).Value = 0.49

# Row 19: id=18, name=листья салата
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "листья салата"
$ws.Cells.Item(19, 3).Value = 10000
$ws.Cells.Item(19, 4).Value = "г"
$ws.Cells.Item(19, 5).Value = 0.3

# Row 20: id=19, name=лосось
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "лосось"
$ws.Cells.Item(20, 3).Value = 10000
$ws.Cells.Item(20, 4).Value = "г"
$ws.Cells.Item(20, 5).Value = 2

# Row 21: id=20, name=майонез
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "майонез"
$ws.Cells.Item(21, 3).Value = 10000
$ws.Cells.Item(21, 4).Value = "г"
$ws.Cells.Item(21, 5).Value = 0.17

# Row 22: id=21, name=масага кр
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "масага кр"
$ws.Cells.Item(22, 3).Value = 10000
$ws.Cells.Item(22, 4).Value = "г"
$ws.Cells.Item(22, 5).Value = 1.08

# Row 23: id=22, name=мицукан
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "мицукан"
$ws.Cells.Item(23, 3).Value = 10000
$ws.Cells.Item(23, 4).Value = "г"
$ws.Cells.Item(23, 5).Value = 0.12

# Row 24: id=23, name=мука
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "мука"
$ws.Cells.Item(24, 3).Value = 10000
$ws.Cells.Item(24, 4).Value = "г"
$ws.Cells.Item(24, 5).Value = 0.23

# Row 25: id=24, name=нори
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "нори"
$ws.Cells.Item(25, 3).Value = 10000
$ws.Cells.Item(25, 4).Value = "г"
$ws.Cells.Item(25, 5).Value = 1.1

# Row 26: id=25, name=огурец
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "огурец"
$ws.Cells.Item(26, 3).Value = 10000
$ws.Cells.Item(26, 4).Value = "г"
$ws.Cells.Item(26, 5).Value = 0.07000000000000001

# Row 27: id=26, name=огурцы
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "огурцы"
$ws.Cells.Item(27, 3).Value = 10000
$ws.Cells.Item(27, 4).Value = "г"
$ws.Cells.Item(27, 5).Value = 0.07000000000000001

# Row 28: id=27, name=омлет
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "омлет"
$ws.Cells.Item(28, 3).Value = 10000
$ws.Cells.Item(28, 4).Value = "г"
$ws.Cells.Item(28, 5).Value = 0.25

# Row 29: id=28, name=остр соус
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "остр соус"
$ws.Cells.Item(29, 3).Value = 10000
$ws.Cells.Item(29, 4).Value = "г"
$ws.Cells.Item(29, 5).Value = 0.35

# Row 30: id=29, name=перец
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "перец"
$ws.Cells.Item(30, 3).Value = 10000
$ws.Cells.Item(30, 4).Value = "г"
$ws.Cells.Item(30, 5).Value = 0.13

# Row 31: id=30, name=помидор
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "помидор"
$ws.Cells.Item(31, 3).Value = 10000
$ws.Cells.Item(31, 4).Value = "г"
$ws.Cells.Item(31, 5).Value = 0.05

# Row 32: id=31, name=рис
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "рис"
$ws.Cells.Item(32, 3).Value = 10000
$ws.Cells.Item(32, 4).Value = "г"
$ws.Cells.Item(32, 5).Value = 0.18

# Row 33: id=32, name=сахар
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "сахар"
$ws.Cells.Item(33, 3).Value = 10000
$ws.Cells.Item(33, 4).Value = "г"
$ws.Cells.Item(33, 5).Value = 0.08

# Row 34: id=33, name=семга
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "семга"
$ws.Cells.Item(34, 3).Value = 10000
$ws.Cells.Item(34, 4).Value = "г"
$ws.Cells.Item(34, 5).Value = 2.19

# Row 35: id=34, name=соевый соус
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "соевый соус"
$ws.Cells.Item(35, 3).Value = 10000
$ws.Cells.Item(35, 4).Value = "г"
$ws.Cells.Item(35, 5).Value = 0.35

# Row 36: id=35, name=соль
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "соль"
$ws.Cells.Item(36, 3).Value = 10000
$ws.Cells.Item(36, 4).Value = "г"
$ws.Cells.Item(36, 5).Value = 0.01

# Row 37: id=36, name=сприн тесто
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "сприн тесто"
$ws.Cells.Item(37, 3).Value = 10000
$ws.Cells.Item(37, 4).Value = "г"
$ws.Cells.Item(37, 5).Value = 10

# Row 38: id=37, name=сухари
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "сухари"
$ws.Cells.Item(38, 3).Value = 10000
$ws.Cells.Item(38, 4).Value = "г"
$ws.Cells.Item(38, 5).Value = 0.15

# Row 39: id=38, name=сыр пармизан
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "сыр пармизан"
$ws.Cells.Item(39, 3).Value = 10000
$ws.Cells.Item(39, 4).Value = "г"
$ws.Cells.Item(39, 5).Value = 0.3

# Row 40: id=39, name=сыр соус
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "сыр соус"
$ws.Cells.Item(40, 3).Value = 10000
$ws.Cells.Item(40, 4).Value = "г"
$ws.Cells.Item(40, 5).Value = 0.34

# Row 41: id=40, name=сыр твор
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "сыр твор"
$ws.Cells.Item(41, 3).Value = 10000
$ws.Cells.Item(41, 4).Value = "г"
$ws.Cells.Item(41, 5).Value = 0.8

# Row 42: id=41, name=сыр творож
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "сыр творож"
$ws.Cells.Item(42, 3).Value = 10000
$ws.Cells.Item(42, 4).Value = "г"
$ws.Cells.Item(42, 5).Value = 0.8

# Row 43: id=42, name=сыр товр
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "сыр товр"
$ws.Cells.Item(43, 3).Value = 10000
$ws.Cells.Item(43, 4).Value = "г"
$ws.Cells.Item(43, 5).Value = 0.8

# Row 44: id=43, name=сыр чеддер
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "сыр чеддер"
$ws.Cells.Item(44, 3).Value = 10000
$ws.Cells.Item(44, 4).Value = "г"
$ws.Cells.Item(44, 5).Value = 1.24

# Row 45: id=44, name=сырный соус
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "сырный соус"
$ws.Cells.Item(45, 3).Value = 10000
$ws.Cells.Item(45, 4).Value = "г"
$ws.Cells.Item(45, 5).Value = 0.34

# Row 46: id=45, name=сырный соус 350
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "сырный соус 350"
$ws.Cells.Item(46, 3).Value = 10000
$ws.Cells.Item(46, 4).Value = "г"
$ws.Cells.Item(46, 5).Value = 0.34

# Row 47: id=46, name=твор сыр
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "твор сыр"
$ws.Cells.Item(47, 3).Value = 10000
$ws.Cells.Item(47, 4).Value = "г"
$ws.Cells.Item(47, 5).Value = 0.8

# Row 48: id=47, name=угорь
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "угорь"
$ws.Cells.Item(48, 3).Value = 10000
$ws.Cells.Item(48, 4).Value = "г"
$ws.Cells.Item(48, 5).Value = 1.74

# Row 49: id=48, name=унаги соус
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "унаги соус"
$ws.Cells.Item(49, 3).Value = 10000
$ws.Cells.Item(49, 4).Value = "г"
$ws.Cells.Item(49, 5).Value = 0.47

# Row 50: id=49, name=чипсы
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "чипсы"
$ws.Cells.Item(50, 3).Value = 10000
$ws.Cells.Item(50, 4).Value = "г"
$ws.Cells.Item(50, 5).Value = 1.07

# Row 51: id=50, name=ширачи
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "ширачи"
$ws.Cells.Item(51, 3).Value = 10000
$ws.Cells.Item(51, 4).Value = "г"
$ws.Cells.Item(51, 5).Value = 1.48

# Row 52: id=51, name=шоколад
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "шоколад"
$ws.Cells.Item(52, 3).Value = 10000
$ws.Cells.Item(52, 4).Value = "г"
$ws.Cells.Item(52, 5).Value = 1.79

# Row 53: id=52, name=яблоко
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "яблоко"
$ws.Cells.Item(53, 3).Value = 10000
$ws.Cells.Item(53, 4).Value = "г"
$ws.Cells.Item(53, 5).Value = 0.07000000000000001

# Row 54: id=53, name=яйцо
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "яйцо"
$ws.Cells.Item(54, 3).Value = 10000
$ws.Cells.Item(54, 4).Value = "г"
$ws.Cells.Item(54, 5).Value = 0.13
